$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that look like plain numbers stay as text,
# matching the original data which stores prices as formatted text (e.g. "10.00", "0.130").
$ws.Range("D2").Value = '68.350.98'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '2.549.49'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.34'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.53'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("D9").Value = '2.546.35'
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("E13").Value = '  -2.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.58'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '2.959.11'
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '68.181.31'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.38'
$ws.Range("E18").Value = '  +139.30%  '
$ws.Range("D19").Value = '2.559.08'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.90'
$ws.Range("E20").Value = '  +4.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.08'
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '371.93'
$ws.Range("E22").Value = '  +3.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.59'
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.22'
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("D29").Value = '2.664.66'
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").Value = '0.0₃0972'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '541.52'
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.34'
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.26'
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.33'
$ws.Range("E39").Value = '  +3.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.63'
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.17'
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.79'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.351'
$ws.Range("E43").Value = '  -1.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.46'
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.05'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0282'
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.73'
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.554'
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.73'
$ws.Range("E51").Value = '  +2.05%  '
